$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8, column C: the phone number was stored as text before; it should
# now be stored as a genuine number.
$ws.Range("C8").Value = 5511970603441

# New row 9: Julio / Julio@Julio.com.br / phone number kept as text.
$ws.Range("A9").Value = "Julio"
$ws.Range("B9").Value = "Julio@Julio.com.br"

# Use a leading apostrophe so Excel stores the phone number as text
# (matching the original C column convention), then reset the cell style
# back to Normal so no stray "quote prefix" formatting is left behind.
$ws.Range("C9").Value = "'5511970603441"
$ws.Range("C9").Style = "Normal"
